$d = $word.ActiveDocument

$old = "1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005.2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018.3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012.4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014.5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;"

$new = "1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005.^l2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018.^l3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012.^l4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014.^l5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
